# Fix the "Intercambio_comercial" (v_expo - v_impo) pipeline bug: the sheet
# previously summed exports and imports (v_expo + v_impo) into column T
# instead of taking the trade balance (v_expo - v_impo). Column U (and the
# p_expo_media / p_impo_media helper columns Q and S) are yearly "media"
# (average) columns per calendar month and are recomputed from the corrected
# T (and from the corrected P/R for row 122) after the fix. A missing data row
# for December 2021 (row 133) is also appended by the corrected CSV-concat step.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row number, then a map of column letter -> corrected value.
$rowFixes = @(
  @{ Row = 2; Cols = @{ 'T' = 295.9260519999998; 'U' = -29.72988752999998 } },
  @{ Row = 3; Cols = @{ 'T' = 35.5403780000006; 'U' = -29.72988752999998 } },
  @{ Row = 4; Cols = @{ 'T' = 24.816333000000668; 'U' = -29.72988752999998 } },
  @{ Row = 5; Cols = @{ 'T' = -1373.3178850000004; 'U' = -29.72988752999998 } },
  @{ Row = 6; Cols = @{ 'T' = -566.7125499999997; 'U' = -29.72988752999998 } },
  @{ Row = 7; Cols = @{ 'T' = -242.36941700000034; 'U' = -29.72988752999998 } },
  @{ Row = 8; Cols = @{ 'T' = -53.361288160000186; 'U' = -29.72988752999998 } },
  @{ Row = 9; Cols = @{ 'T' = -933.0235387800003; 'U' = -29.72988752999998 } },
  @{ Row = 10; Cols = @{ 'T' = 372.9856728300001; 'U' = -29.72988752999998 } },
  @{ Row = 11; Cols = @{ 'T' = 1044.03497147; 'U' = -29.72988752999998 } },
  @{ Row = 12; Cols = @{ 'T' = 1068.45250881; 'U' = -29.72988752999998 } },
  @{ Row = 13; Cols = @{ 'T' = 633.7390439999999; 'U' = 197.80608694090907 } },
  @{ Row = 14; Cols = @{ 'T' = 1279.5526310000005; 'U' = 197.80608694090907 } },
  @{ Row = 15; Cols = @{ 'T' = -268.8669310000005; 'U' = 197.80608694090907 } },
  @{ Row = 16; Cols = @{ 'T' = -889.7143550000001; 'U' = 197.80608694090907 } },
  @{ Row = 17; Cols = @{ 'T' = -212.77615699999978; 'U' = 197.80608694090907 } },
  @{ Row = 18; Cols = @{ 'T' = 44.7195670000001; 'U' = 197.80608694090907 } },
  @{ Row = 19; Cols = @{ 'T' = -217.82092465000005; 'U' = 197.80608694090907 } },
  @{ Row = 20; Cols = @{ 'T' = -891.7662172099999; 'U' = 197.80608694090907 } },
  @{ Row = 21; Cols = @{ 'T' = 449.7158926099996; 'U' = 197.80608694090907 } },
  @{ Row = 22; Cols = @{ 'T' = 1186.7620058900015; 'U' = 197.80608694090907 } },
  @{ Row = 23; Cols = @{ 'T' = 1062.3224007099984; 'U' = 197.80608694090907 } },
  @{ Row = 24; Cols = @{ 'T' = 456.52445700000044; 'U' = 190.671605229091 } },
  @{ Row = 25; Cols = @{ 'T' = 1053.469376; 'U' = 190.671605229091 } },
  @{ Row = 26; Cols = @{ 'T' = 540.7792529999997; 'U' = 190.671605229091 } },
  @{ Row = 27; Cols = @{ 'T' = -476.9265720000003; 'U' = 190.671605229091 } },
  @{ Row = 28; Cols = @{ 'T' = -677.3522430000003; 'U' = 190.671605229091 } },
  @{ Row = 29; Cols = @{ 'T' = -135.47755499999948; 'U' = 190.671605229091 } },
  @{ Row = 30; Cols = @{ 'T' = -910.3121932599997; 'U' = 190.671605229091 } },
  @{ Row = 31; Cols = @{ 'T' = -562.2848089600002; 'U' = 190.671605229091 } },
  @{ Row = 32; Cols = @{ 'T' = 1180.6592827300005; 'U' = 190.671605229091 } },
  @{ Row = 33; Cols = @{ 'T' = 1228.48143192; 'U' = 190.671605229091 } },
  @{ Row = 34; Cols = @{ 'T' = 399.8272290900004; 'U' = 190.671605229091 } },
  @{ Row = 35; Cols = @{ 'T' = 1509.9788760000001; 'U' = 782.993362858182 } },
  @{ Row = 36; Cols = @{ 'T' = 1945.1494040000007; 'U' = 782.993362858182 } },
  @{ Row = 37; Cols = @{ 'T' = 701.9412199999997; 'U' = 782.993362858182 } },
  @{ Row = 38; Cols = @{ 'T' = 822.3911199999993; 'U' = 782.993362858182 } },
  @{ Row = 39; Cols = @{ 'T' = 222.33081000000038; 'U' = 782.993362858182 } },
  @{ Row = 40; Cols = @{ 'T' = 316.90995400000065; 'U' = 782.993362858182 } },
  @{ Row = 41; Cols = @{ 'T' = -106.47264349000034; 'U' = 782.993362858182 } },
  @{ Row = 42; Cols = @{ 'T' = -889.2541233899992; 'U' = 782.993362858182 } },
  @{ Row = 43; Cols = @{ 'T' = 1164.8938146300006; 'U' = 782.993362858182 } },
  @{ Row = 44; Cols = @{ 'T' = 1454.8669904399999; 'U' = 782.993362858182 } },
  @{ Row = 45; Cols = @{ 'T' = 1470.1915692499997; 'U' = 782.993362858182 } },
  @{ Row = 46; Cols = @{ 'T' = 1443.0515809999997; 'U' = 831.1230008209088 } },
  @{ Row = 47; Cols = @{ 'T' = 1093.7342979999994; 'U' = 831.1230008209088 } },
  @{ Row = 48; Cols = @{ 'T' = 1365.164707; 'U' = 831.1230008209088 } },
  @{ Row = 49; Cols = @{ 'T' = 1401.6069390000002; 'U' = 831.1230008209088 } },
  @{ Row = 50; Cols = @{ 'T' = 202.60902799999985; 'U' = 831.1230008209088 } },
  @{ Row = 51; Cols = @{ 'T' = 509.85932499999944; 'U' = 831.1230008209088 } },
  @{ Row = 52; Cols = @{ 'T' = -571.1389577699993; 'U' = 831.1230008209088 } },
  @{ Row = 53; Cols = @{ 'T' = -1284.94533607; 'U' = 831.1230008209088 } },
  @{ Row = 54; Cols = @{ 'T' = 1398.9011870599998; 'U' = 831.1230008209088 } },
  @{ Row = 55; Cols = @{ 'T' = 1911.8206798799993; 'U' = 831.1230008209088 } },
  @{ Row = 56; Cols = @{ 'T' = 1671.689557929999; 'U' = 831.1230008209088 } },
  @{ Row = 57; Cols = @{ 'T' = 1112.4346510000005; 'U' = 662.1667342636365 } },
  @{ Row = 58; Cols = @{ 'T' = 788.6734890000007; 'U' = 662.1667342636365 } },
  @{ Row = 59; Cols = @{ 'T' = 840.5952729999999; 'U' = 662.1667342636365 } },
  @{ Row = 60; Cols = @{ 'T' = 1355.9314939999995; 'U' = 662.1667342636365 } },
  @{ Row = 61; Cols = @{ 'T' = 339.4167610000004; 'U' = 662.1667342636365 } },
  @{ Row = 62; Cols = @{ 'T' = 258.5719989999998; 'U' = 662.1667342636365 } },
  @{ Row = 63; Cols = @{ 'T' = -739.6614130399994; 'U' = 662.1667342636365 } },
  @{ Row = 64; Cols = @{ 'T' = -327.593278620001; 'U' = 662.1667342636365 } },
  @{ Row = 65; Cols = @{ 'T' = 1067.7531275500005; 'U' = 662.1667342636365 } },
  @{ Row = 66; Cols = @{ 'T' = 1486.6397080700008; 'U' = 662.1667342636365 } },
  @{ Row = 67; Cols = @{ 'T' = 1101.0722659399999; 'U' = 662.1667342636365 } },
  @{ Row = 68; Cols = @{ 'T' = 367.06863599999997; 'U' = 400.44510731909105 } },
  @{ Row = 69; Cols = @{ 'T' = 1310.2453569999998; 'U' = 400.44510731909105 } },
  @{ Row = 70; Cols = @{ 'T' = -57.489841000000524; 'U' = 400.44510731909105 } },
  @{ Row = 71; Cols = @{ 'T' = 271.83666300000004; 'U' = 400.44510731909105 } },
  @{ Row = 72; Cols = @{ 'T' = -231.67406799999935; 'U' = 400.44510731909105 } },
  @{ Row = 73; Cols = @{ 'T' = 283.9796700000006; 'U' = 400.44510731909105 } },
  @{ Row = 74; Cols = @{ 'T' = -738.6524935699999; 'U' = 400.44510731909105 } },
  @{ Row = 75; Cols = @{ 'T' = -767.5699735600001; 'U' = 400.44510731909105 } },
  @{ Row = 76; Cols = @{ 'T' = 950.4186272499992; 'U' = 400.44510731909105 } },
  @{ Row = 77; Cols = @{ 'T' = 1480.2650479000004; 'U' = 400.44510731909105 } },
  @{ Row = 78; Cols = @{ 'T' = 1536.4685554900016; 'U' = 400.44510731909105 } },
  @{ Row = 79; Cols = @{ 'T' = 657.4109950000002; 'U' = 514.2631408127269 } },
  @{ Row = 80; Cols = @{ 'T' = 1174.6786240000001; 'U' = 514.2631408127269 } },
  @{ Row = 81; Cols = @{ 'T' = 125.64711799999986; 'U' = 514.2631408127269 } },
  @{ Row = 82; Cols = @{ 'T' = 775.5758309999992; 'U' = 514.2631408127269 } },
  @{ Row = 83; Cols = @{ 'T' = -469.2139590000006; 'U' = 514.2631408127269 } },
  @{ Row = 84; Cols = @{ 'T' = 590.0863879999997; 'U' = 514.2631408127269 } },
  @{ Row = 85; Cols = @{ 'T' = -1044.3547655700022; 'U' = 514.2631408127269 } },
  @{ Row = 86; Cols = @{ 'T' = -1112.4908280700001; 'U' = 514.2631408127269 } },
  @{ Row = 87; Cols = @{ 'T' = 1167.9414548000004; 'U' = 514.2631408127269 } },
  @{ Row = 88; Cols = @{ 'T' = 1446.4014471699998; 'U' = 514.2631408127269 } },
  @{ Row = 89; Cols = @{ 'T' = 2345.212243609999; 'U' = 514.2631408127269 } },
  @{ Row = 90; Cols = @{ 'T' = 829.6939029999994; 'U' = 543.7406717445454 } },
  @{ Row = 91; Cols = @{ 'T' = 948.5082680000005; 'U' = 543.7406717445454 } },
  @{ Row = 92; Cols = @{ 'T' = 515.8838089999999; 'U' = 543.7406717445454 } },
  @{ Row = 93; Cols = @{ 'T' = 41.1527759999999; 'U' = 543.7406717445454 } },
  @{ Row = 94; Cols = @{ 'T' = -325.07707800000026; 'U' = 543.7406717445454 } },
  @{ Row = 95; Cols = @{ 'T' = 308.1335740000004; 'U' = 543.7406717445454 } },
  @{ Row = 96; Cols = @{ 'T' = -700.1108071900007; 'U' = 543.7406717445454 } },
  @{ Row = 97; Cols = @{ 'T' = 336.85376734000056; 'U' = 543.7406717445454 } },
  @{ Row = 98; Cols = @{ 'T' = 1744.0952103099999; 'U' = 543.7406717445454 } },
  @{ Row = 99; Cols = @{ 'T' = 597.5067276999998; 'U' = 543.7406717445454 } },
  @{ Row = 100; Cols = @{ 'T' = 1684.5072390299993; 'U' = 543.7406717445454 } },
  @{ Row = 101; Cols = @{ 'T' = 1132.4116620000004; 'U' = 421.6564482963638 } },
  @{ Row = 102; Cols = @{ 'T' = 632.6058860000003; 'U' = 421.6564482963638 } },
  @{ Row = 103; Cols = @{ 'T' = -803.9730509999999; 'U' = 421.6564482963638 } },
  @{ Row = 104; Cols = @{ 'T' = 299.76638399999956; 'U' = 421.6564482963638 } },
  @{ Row = 105; Cols = @{ 'T' = -25.03436800000054; 'U' = 421.6564482963638 } },
  @{ Row = 106; Cols = @{ 'T' = -46.2107019999994; 'U' = 421.6564482963638 } },
  @{ Row = 107; Cols = @{ 'T' = -928.6008027199996; 'U' = 421.6564482963638 } },
  @{ Row = 108; Cols = @{ 'T' = 324.6519118900005; 'U' = 421.6564482963638 } },
  @{ Row = 109; Cols = @{ 'T' = 1767.6625532400003; 'U' = 421.6564482963638 } },
  @{ Row = 110; Cols = @{ 'T' = 669.6034968400008; 'U' = 421.6564482963638 } },
  @{ Row = 111; Cols = @{ 'T' = 1615.3379610099992; 'U' = 421.6564482963638 } },
  @{ Row = 112; Cols = @{ 'Q' = 5323.552942100909; 'S' = 5089.209807263636; 'T' = 263.41857600000003; 'U' = 1281.0704075645453 } },
  @{ Row = 113; Cols = @{ 'Q' = 5323.552942100909; 'S' = 5089.209807263636; 'T' = 869.9641939999992; 'U' = 1281.0704075645453 } },
  @{ Row = 114; Cols = @{ 'Q' = 5323.552942100909; 'S' = 5089.209807263636; 'T' = -1002.0742579999996; 'U' = 1281.0704075645453 } },
  @{ Row = 115; Cols = @{ 'Q' = 5323.552942100909; 'S' = 5089.209807263636; 'T' = 238.20841500000006; 'U' = 1281.0704075645453 } },
  @{ Row = 116; Cols = @{ 'Q' = 5323.552942100909; 'S' = 5089.209807263636; 'T' = -700.9279040000001; 'U' = 1281.0704075645453 } },
  @{ Row = 117; Cols = @{ 'Q' = 5323.552942100909; 'S' = 5089.209807263636; 'T' = 114.27996600000006; 'U' = 1281.0704075645453 } },
  @{ Row = 118; Cols = @{ 'Q' = 5323.552942100909; 'S' = 5089.209807263636; 'T' = -1495.3551513000011; 'U' = 1281.0704075645453 } },
  @{ Row = 119; Cols = @{ 'Q' = 5323.552942100909; 'S' = 5089.209807263636; 'T' = 986.5954697800007; 'U' = 1281.0704075645453 } },
  @{ Row = 120; Cols = @{ 'Q' = 5323.552942100909; 'S' = 5089.209807263636; 'T' = 2483.98902806; 'U' = 1281.0704075645453 } },
  @{ Row = 121; Cols = @{ 'Q' = 5323.552942100909; 'S' = 5089.209807263636; 'T' = 385.6761476699994; 'U' = 1281.0704075645453 } },
  @{ Row = 122; Cols = @{ 'P' = 6191.0; 'Q' = 5323.552942100909; 'R' = 5757.0; 'S' = 5089.209807263636; 'T' = 11948.0; 'U' = 1281.0704075645453 } },
  @{ Row = 123; Cols = @{ 'Q' = 5015.8882719172725; 'S' = 4746.862200523636; 'T' = 318.867201; 'U' = 1393.4806168481819 } },
  @{ Row = 124; Cols = @{ 'Q' = 5015.8882719172725; 'S' = 4746.862200523636; 'T' = 876.14732; 'U' = 1393.4806168481819 } },
  @{ Row = 125; Cols = @{ 'Q' = 5015.8882719172725; 'S' = 4746.862200523636; 'T' = -461.14007200000015; 'U' = 1393.4806168481819 } },
  @{ Row = 126; Cols = @{ 'Q' = 5015.8882719172725; 'S' = 4746.862200523636; 'T' = 201.8685800000003; 'U' = 1393.4806168481819 } },
  @{ Row = 127; Cols = @{ 'Q' = 5015.8882719172725; 'S' = 4746.862200523636; 'T' = -974.56072; 'U' = 1393.4806168481819 } },
  @{ Row = 128; Cols = @{ 'Q' = 5015.8882719172725; 'S' = 4746.862200523636; 'T' = 54.665218000000095; 'U' = 1393.4806168481819 } },
  @{ Row = 129; Cols = @{ 'Q' = 5015.8882719172725; 'S' = 4746.862200523636; 'T' = -787.5095916700011; 'U' = 1393.4806168481819 } },
  @{ Row = 130; Cols = @{ 'Q' = 5015.8882719172725; 'S' = 4746.862200523636; 'T' = 1419.5411408900004; 'U' = 1393.4806168481819 } },
  @{ Row = 131; Cols = @{ 'Q' = 5015.8882719172725; 'S' = 4746.862200523636; 'T' = 2241.2807965300003; 'U' = 1393.4806168481819 } },
  @{ Row = 132; Cols = @{ 'Q' = 5015.8882719172725; 'S' = 4746.862200523636; 'T' = -363.87308742000005; 'U' = 1393.4806168481819 } },
)

$colIndex = @{
  'A' = 1
  'B' = 2
  'C' = 3
  'D' = 4
  'E' = 5
  'F' = 6
  'G' = 7
  'H' = 8
  'I' = 9
  'J' = 10
  'K' = 11
  'L' = 12
  'M' = 13
  'N' = 14
  'O' = 15
  'P' = 16
  'Q' = 17
  'R' = 18
  'S' = 19
  'T' = 20
  'U' = 21
}

foreach ($fix in $rowFixes) {
    $r = $fix.Row
    foreach ($col in $fix.Cols.Keys) {
        $c = $colIndex[$col]
        $ws.Cells.Item($r, $c).Value = $fix.Cols[$col]
    }
}

# Append the new row 133 (Mes_num=12 "dic", Año=2021) produced by the fixed
# dataframe-concatenation step.
$ws.Cells.Item(133, $colIndex["A"]).Value = 12
$ws.Cells.Item(133, $colIndex["B"]).Value = "dic"
$ws.Cells.Item(133, $colIndex["C"]).Value = 2021
$ws.Cells.Item(133, $colIndex["J"]).Value = 170.00381704316
$ws.Cells.Item(133, $colIndex["K"]).Value = 165.47068406624
$ws.Cells.Item(133, $colIndex["L"]).Value = 102.3848577974986
$ws.Cells.Item(133, $colIndex["M"]).Value = 248.38308174534
$ws.Cells.Item(133, $colIndex["N"]).Value = 123.94421688775
$ws.Cells.Item(133, $colIndex["O"]).Value = 199.6054751216285
$ws.Cells.Item(133, $colIndex["P"]).Value = 6191.0
$ws.Cells.Item(133, $colIndex["Q"]).Value = 5015.8882719172725
$ws.Cells.Item(133, $colIndex["R"]).Value = 5757.0
$ws.Cells.Item(133, $colIndex["S"]).Value = 4746.862200523636
$ws.Cells.Item(133, $colIndex["T"]).Value = 12803.0
$ws.Cells.Item(133, $colIndex["U"]).Value = 1393.4806168481819
